$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the Correspond Handoff/Handback datetimes for the
# d46f061f... row (row 3)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-14 08:33:16"
$wsZh.Range("H3").Value = "2016-03-14 08:33:33"

# de-de sheet: update the Correspond Handoff/Handback datetimes for the
# d46f061f... row (row 3)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-14 08:33:19"
$wsDe.Range("H3").Value = "2016-03-14 08:33:38"
